$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column L (rows 3-11) into new column M, replicating formatting + values
$ws.Range("L3:L11").Copy($ws.Range("M3:M11")) | Out-Null

# Set the new year value for the header row
$ws.Range("M4").Value = 2022

# Update selection to N6
$ws.Range("N6").Select() | Out-Null
